# Adding the RES Hourly Production Forecast to the Portfolio
# - Shift all interval timestamps in column A from 29.08.2024 to 24.09.2024 (+26 days)
# - Update the Lookup strings in column D to match the new date
# - Refresh the Prediction values in column C (rows 31-82) with the new forecast

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dayShift = 26

# Shift the date/time serials in column A (rows 2..96) and refresh the
# "Lookup" text in column D to reflect the new date (29.08.2024 -> 24.09.2024).
for ($r = 2; $r -le 96; $r++) {
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value2 = $dateCell.Value2 + $dayShift

    $lookupCell = $ws.Cells.Item($r, 4)
    $lookupCell.Value2 = $lookupCell.Value2.Replace("29.08.2024", "24.09.2024")
}

# Updated Prediction values for rows 31-82 (interval 30-81), the rest of the
# rows keep their existing (zero) values.
$newPredictions = @{
    31 = 0
    32 = 0.022
    33 = 0.059
    34 = 0.104
    35 = 0.156
    36 = 0.212
    37 = 0.267
    38 = 0.322
    39 = 0.374
    40 = 0.424
    41 = 0.47
    42 = 0.511
    43 = 0.548
    44 = 0.582
    45 = 0.612
    46 = 0.638
    47 = 0.663
    48 = 0.6840000000000001
    49 = 0.705
    50 = 0.726
    51 = 0.739
    52 = 0.747
    53 = 0.751
    54 = 0.751
    55 = 0.751
    56 = 0.752
    57 = 0.738
    58 = 0.729
    59 = 0.714
    60 = 0.697
    61 = 0.657
    62 = 0.629
    63 = 0.585
    64 = 0.537
    65 = 0.484
    66 = 0.457
    67 = 0.426
    68 = 0.381
    69 = 0.306
    70 = 0.271
    71 = 0.232
    72 = 0.185
    73 = 0.153
    74 = 0.125
    75 = 0.091
    76 = 0.065
    77 = 0.05
    78 = 0.04
    79 = 0
    80 = 0
    81 = 0
    82 = 0
}

foreach ($row in $newPredictions.Keys) {
    $ws.Cells.Item($row, 3).Value2 = $newPredictions[$row]
}
